$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the constant-column values for the brand new row 229 -------------
# Every data row (2..228) shares identical values in these columns, so the
# newly appended row copies them straight from row 228.
$ws.Range("A229").Value2 = $ws.Range("A228").Value2
$ws.Range("B229").Value2 = $ws.Range("B228").Value2
$ws.Range("C229").Value2 = $ws.Range("C228").Value2
$ws.Range("E229").Value2 = $ws.Range("E228").Value2
$ws.Range("F229").Value2 = $ws.Range("F228").Value2
$ws.Range("G229").Value2 = $ws.Range("G228").Value2
$ws.Range("H229").Value2 = $ws.Range("H228").Value2
$ws.Range("I229").Value2 = $ws.Range("I228").Value2
$ws.Range("N229").Value2 = $ws.Range("N228").Value2
$ws.Range("O229").Value2 = $ws.Range("O228").Value2
$ws.Range("Q229").Value2 = $ws.Range("Q228").Value2
$ws.Range("R229").Value2 = $ws.Range("R228").Value2

# --- Shift the weekly readings (Fecha / Volumen / Precios) down one row ---
# A new weekly record is inserted at row 29 (date 44550); every row from the
# old row 29 through the old row 228 slides down to make room, with the
# former last row (228) landing in the brand-new row 229.

$oldCount = 200              # old rows 29..228 inclusive
$newCount = 201               # new rows 29..229 inclusive

$dOld = $ws.Range("D29:D228").Value2
$jmOld = $ws.Range("J29:M228").Value2
$pOld = $ws.Range("P29:P228").Value2

$dNew = New-Object 'object[,]' $newCount,1
$jmNew = New-Object 'object[,]' $newCount,4
$pNew = New-Object 'object[,]' $newCount,1

# Row 29 (array index 0) keeps its own Volumen/Precio figures, only the
# date changes to the newly recorded week.
$dNew[0,0] = 44550
$jmNew[0,0] = $jmOld[1,1]
$jmNew[0,1] = $jmOld[1,2]
$jmNew[0,2] = $jmOld[1,3]
$jmNew[0,3] = $jmOld[1,4]
$pNew[0,0] = $pOld[1,1]

# Array indices 1..oldCount (0-based dNew/jmNew/pNew) receive the old rows
# 29..228 (1-based dOld/jmOld/pOld indices 1..oldCount) in order, i.e. a
# one-row-down shift.
for ($i = 1; $i -le $oldCount; $i++) {
    $dNew[$i,0] = $dOld[$i,1]
    $jmNew[$i,0] = $jmOld[$i,1]
    $jmNew[$i,1] = $jmOld[$i,2]
    $jmNew[$i,2] = $jmOld[$i,3]
    $jmNew[$i,3] = $jmOld[$i,4]
    $pNew[$i,0] = $pOld[$i,1]
}

$ws.Range("D29:D229").Value2 = $dNew
$ws.Range("J29:M229").Value2 = $jmNew
$ws.Range("P29:P229").Value2 = $pNew

# The new row 229 needs the same date formatting as every other Fecha cell.
$ws.Range("D229").NumberFormat = $ws.Range("D228").NumberFormat
